$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B4").Value = "Account Number"
$ws.Range("C4").Value = "Amount"
$ws.Range("B4:C4").Font.Bold = $true

# Account numbers (column B) and amounts (column C)
$accounts = @("A0000", "A1111", "A2222", "A3333", "A4444", "A5555", "A7777", "A8888", "A9999")
$amounts  = @(9836, 0, 8508, 3349, 0, 0, 5926, 8472, 0)

for ($i = 0; $i -lt $accounts.Length; $i++) {
    $row = 5 + $i
    $cell = $ws.Cells.Item($row, 2)
    if ($i -lt 2) {
        $cell.Value = "'" + $accounts[$i]
    } else {
        $cell.Value = $accounts[$i]
    }
    $cell.HorizontalAlignment = -4152
    $ws.Cells.Item($row, 3).Value = $amounts[$i]
}

$ws.Columns.Item(2).ColumnWidth = 16
$ws.Columns.Item(3).ColumnWidth = 16

$ws.Range("B4").Select()
